$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9142543333333334
$ws.Range("H2").Value = 2.742763
$ws.Range("I2").Value = 0.1175834869881751
$ws.Range("J2").Value = 0.1175834869881751
$ws.Range("M2").Value = 255.0443116666667
$ws.Range("N2").Value = 765.132935
$ws.Range("O2").Value = 0.863617428561108
$ws.Range("P2").Value = 0.8636174285611079
$ws.Range("Q2").Value = 233.1753671332672
$ws.Range("R2").Value = 2098.578304199405
$ws.Range("S2").Value = 0.1015471486739763
$ws.Range("T2").Value = 0.1015471486739762

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9142543333333334
$ws.Range("H3").Value = 2.742763
$ws.Range("I3").Value = 0.1175834869881751
$ws.Range("J3").Value = 0.1175834869881751
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("O3").Value = 0.003031431940796009
$ws.Range("P3").Value = 0.003031431940796009
$ws.Range("Q3").Value = 0.8184819253964446
$ws.Range("R3").Value = 7.366337328568001
$ws.Range("S3").Value = 0.0003564463381661258
$ws.Range("T3").Value = 0.0003564463381661258

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9142543333333334
$ws.Range("H4").Value = 2.742763
$ws.Range("I4").Value = 0.1175834869881751
$ws.Range("J4").Value = 0.1175834869881751
$ws.Range("M4").Value = 7.050555333333333
$ws.Range("N4").Value = 21.151666
$ws.Range("O4").Value = 0.02387421396349043
$ws.Range("P4").Value = 0.02387421396349043
$ws.Range("Q4").Value = 6.446000765906444
$ws.Range("R4").Value = 58.014006893158
$ws.Range("S4").Value = 0.002807213326928985
$ws.Range("T4").Value = 0.002807213326928984

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9142543333333334
$ws.Range("H5").Value = 2.742763
$ws.Range("I5").Value = 0.1175834869881751
$ws.Range("J5").Value = 0.1175834869881751
$ws.Range("M5").Value = 32.33082866666666
$ws.Range("N5").Value = 96.99248599999999
$ws.Range("O5").Value = 0.1094769255346056
$ws.Range("P5").Value = 0.1094769255346056
$ws.Range("Q5").Value = 29.55860020875755
$ws.Range("R5").Value = 266.027401878818
$ws.Range("S5").Value = 0.01287267864910371
$ws.Range("T5").Value = 0.01287267864910371

$ws.Range("I6").Value = 0.6206849497708361
$ws.Range("J6").Value = 0.620684949770836
$ws.Range("M6").Value = 255.0443116666667
$ws.Range("N6").Value = 765.132935
$ws.Range("O6").Value = 0.863617428561108
$ws.Range("P6").Value = 0.8636174285611079
$ws.Range("Q6").Value = 1230.856855363228
$ws.Range("R6").Value = 11077.71169826905
$ws.Range("S6").Value = 0.5360343402676699
$ws.Range("T6").Value = 0.5360343402676698

$ws.Range("I7").Value = 0.6206849497708361
$ws.Range("J7").Value = 0.620684949770836
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("O7").Value = 0.003031431940796009
$ws.Range("P7").Value = 0.003031431940796009
$ws.Range("Q7").Value = 4.320499636178667
$ws.Range("S7").Value = 0.001881564181906679
$ws.Range("T7").Value = 0.001881564181906678

$ws.Range("I8").Value = 0.6206849497708361
$ws.Range("J8").Value = 0.620684949770836
$ws.Range("M8").Value = 7.050555333333333
$ws.Range("N8").Value = 21.151666
$ws.Range("O8").Value = 0.02387421396349043
$ws.Range("P8").Value = 0.02387421396349043
$ws.Range("Q8").Value = 34.02633961698866
$ws.Range("R8").Value = 306.237056552898
$ws.Range("S8").Value = 0.01481836529474725
$ws.Range("T8").Value = 0.01481836529474725

$ws.Range("I9").Value = 0.6206849497708361
$ws.Range("J9").Value = 0.620684949770836
$ws.Range("M9").Value = 32.33082866666666
$ws.Range("N9").Value = 96.99248599999999
$ws.Range("O9").Value = 0.1094769255346056
$ws.Range("P9").Value = 0.1094769255346056
$ws.Range("Q9").Value = 156.0302280175953
$ws.Range("R9").Value = 1404.272052158358
$ws.Range("S9").Value = 0.06795068002651226
$ws.Range("T9").Value = 0.06795068002651225

$ws.Range("G10").Value = 1.986145
$ws.Range("H10").Value = 5.958435
$ws.Range("I10").Value = 0.2554407961214246
$ws.Range("J10").Value = 0.2554407961214246
$ws.Range("M10").Value = 255.0443116666667
$ws.Range("N10").Value = 765.132935
$ws.Range("O10").Value = 0.863617428561108
$ws.Range("P10").Value = 0.8636174285611079
$ws.Range("Q10").Value = 506.5549843951916
$ws.Range("R10").Value = 4558.994859556725
$ws.Range("S10").Value = 0.2206031234959869
$ws.Range("T10").Value = 0.2206031234959869

$ws.Range("G11").Value = 1.986145
$ws.Range("H11").Value = 5.958435
$ws.Range("I11").Value = 0.2554407961214246
$ws.Range("J11").Value = 0.2554407961214246
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("O11").Value = 0.003031431940796009
$ws.Range("P11").Value = 0.003031431940796009
$ws.Range("Q11").Value = 1.778087042573333
$ws.Range("R11").Value = 16.00278338316
$ws.Range("S11").Value = 0.0007743513883448477
$ws.Range("T11").Value = 0.0007743513883448477

$ws.Range("G12").Value = 1.986145
$ws.Range("H12").Value = 5.958435
$ws.Range("I12").Value = 0.2554407961214246
$ws.Range("J12").Value = 0.2554407961214246
$ws.Range("M12").Value = 7.050555333333333
$ws.Range("N12").Value = 21.151666
$ws.Range("O12").Value = 0.02387421396349043
$ws.Range("P12").Value = 0.02387421396349043
$ws.Range("Q12").Value = 14.00342522252333
$ws.Range("R12").Value = 126.03082700271
$ws.Range("S12").Value = 0.006098448221607227
$ws.Range("T12").Value = 0.006098448221607226

$ws.Range("G13").Value = 1.986145
$ws.Range("H13").Value = 5.958435
$ws.Range("I13").Value = 0.2554407961214246
$ws.Range("J13").Value = 0.2554407961214246
$ws.Range("M13").Value = 32.33082866666666
$ws.Range("N13").Value = 96.99248599999999
$ws.Range("O13").Value = 0.1094769255346056
$ws.Range("P13").Value = 0.1094769255346056
$ws.Range("Q13").Value = 64.21371370215665
$ws.Range("R13").Value = 577.9234233194098
$ws.Range("S13").Value = 0.02796487301548558
$ws.Range("T13").Value = 0.02796487301548558

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.04891300000000001
$ws.Range("H14").Value = 0.146739
$ws.Range("I14").Value = 0.006290767119564404
$ws.Range("J14").Value = 0.006290767119564403
$ws.Range("M14").Value = 255.0443116666667
$ws.Range("N14").Value = 765.132935
$ws.Range("O14").Value = 0.863617428561108
$ws.Range("P14").Value = 0.8636174285611079
$ws.Range("Q14").Value = 12.47498241655167
$ws.Range("R14").Value = 112.274841748965
$ws.Range("S14").Value = 0.005432816123474979
$ws.Range("T14").Value = 0.005432816123474977

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.04891300000000001
$ws.Range("H15").Value = 0.146739
$ws.Range("I15").Value = 0.006290767119564404
$ws.Range("J15").Value = 0.006290767119564403
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("O15").Value = 0.003031431940796009
$ws.Range("P15").Value = 0.003031431940796009
$ws.Range("Q15").Value = 0.04378913498933334
$ws.Range("R15").Value = 0.394102214904
$ws.Range("S15").Value = 0.00001907003237835684
$ws.Range("T15").Value = 0.00001907003237835683

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.04891300000000001
$ws.Range("H16").Value = 0.146739
$ws.Range("I16").Value = 0.006290767119564404
$ws.Range("J16").Value = 0.006290767119564403
$ws.Range("M16").Value = 7.050555333333333
$ws.Range("N16").Value = 21.151666
$ws.Range("O16").Value = 0.02387421396349043
$ws.Range("P16").Value = 0.02387421396349043
$ws.Range("Q16").Value = 0.3448638130193333
$ws.Range("R16").Value = 3.103774317174
$ws.Range("S16").Value = 0.000150187120206971
$ws.Range("T16").Value = 0.0001501871202069709

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.04891300000000001
$ws.Range("H17").Value = 0.146739
$ws.Range("I17").Value = 0.006290767119564404
$ws.Range("J17").Value = 0.006290767119564403
$ws.Range("M17").Value = 32.33082866666666
$ws.Range("N17").Value = 96.99248599999999
$ws.Range("O17").Value = 0.1094769255346056
$ws.Range("P17").Value = 0.1094769255346056
$ws.Range("Q17").Value = 1.581397822572667
$ws.Range("R17").Value = 14.232580403154
$ws.Range("S17").Value = 0.0006886938435040978
$ws.Range("T17").Value = 0.0006886938435040976
